$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (shifts old row 3 down to row 4),
# so that the new "ECs" row data structure mirrors the diff where
# row 2 is ECs, row 3 is FAPs (formerly row 2 content updated), row 4 is sCs (old row 3 content updated)
$ws.Rows.Item(3).Insert()

# Row 2: ECs -> Ccl11 -> Ackr2 -> FAPs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ccl11"
$ws.Cells.Item(2, 3).Value = "Ackr2"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.6658270000000001
$ws.Cells.Item(2, 8).Value = 1.997481
$ws.Cells.Item(2, 9).Value = 0.002619014661824084
$ws.Cells.Item(2, 10).Value = 0.002619014661824084
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 11.51723133333333
$ws.Cells.Item(2, 14).Value = 34.551694
$ws.Cells.Item(2, 15).Value = 1
$ws.Cells.Item(2, 16).Value = 1
$ws.Cells.Item(2, 17).Value = 7.668483586979334
$ws.Cells.Item(2, 18).Value = 69.016352282814
$ws.Cells.Item(2, 19).Value = 0.002619014661824084
$ws.Cells.Item(2, 20).Value = 0.002619014661824084

# Row 3: FAPs -> Ccl11 -> Ackr2 -> FAPs
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Ccl11"
$ws.Cells.Item(3, 3).Value = "Ackr2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 245.845932
$ws.Cells.Item(3, 8).Value = 737.537796
$ws.Cells.Item(3, 9).Value = 0.967029123868222
$ws.Cells.Item(3, 10).Value = 0.9670291238682222
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 11.51723133333333
$ws.Cells.Item(3, 14).Value = 34.551694
$ws.Cells.Item(3, 15).Value = 1
$ws.Cells.Item(3, 16).Value = 1
$ws.Cells.Item(3, 17).Value = 2831.464471202936
$ws.Cells.Item(3, 18).Value = 25483.18024082642
$ws.Cells.Item(3, 19).Value = 0.967029123868222
$ws.Cells.Item(3, 20).Value = 0.9670291238682222

# Row 4: sCs -> Ccl11 -> Ackr2 -> FAPs
$ws.Cells.Item(4, 1).Value = "sCs"
$ws.Cells.Item(4, 2).Value = "Ccl11"
$ws.Cells.Item(4, 3).Value = "Ackr2"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 7.716294666666666
$ws.Cells.Item(4, 8).Value = 23.148884
$ws.Cells.Item(4, 9).Value = 0.03035186146995388
$ws.Cells.Item(4, 10).Value = 0.03035186146995388
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 11.51723133333333
$ws.Cells.Item(4, 14).Value = 34.551694
$ws.Cells.Item(4, 15).Value = 1
$ws.Cells.Item(4, 16).Value = 1
$ws.Cells.Item(4, 17).Value = 88.87035071216621
$ws.Cells.Item(4, 18).Value = 799.8331564094959
$ws.Cells.Item(4, 19).Value = 0.03035186146995388
$ws.Cells.Item(4, 20).Value = 0.03035186146995388

$wb.Save()
